$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.530.79'
$ws.Range('E2').Value = '  -2.70%  '
$ws.Range('D3').Value = '2.217.11'
$ws.Range('E3').Value = '  -6.50%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '298.06'
$ws.Range('E5').Value = '  -4.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '83.08'
$ws.Range('E6').Value = '  -4.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.513'
$ws.Range('E7').Value = '  -3.51%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.470'
$ws.Range('E9').Value = '  -4.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0778'
$ws.Range('E10').Value = '  -6.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '29.29'
$ws.Range('E11').Value = '  -3.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.65'
$ws.Range('E12').Value = '  -9.80%  '
$ws.Range('E13').Value = '  -2.39%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.32'
$ws.Range('E14').Value = '  -3.10%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '2.556.16'
$ws.Range('E15').Value = '  -6.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.14'
$ws.Range('E16').Value = '  -5.17%  '
$ws.Range('D17').Value = '2.205.11'
$ws.Range('E17').Value = '  -7.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.716'
$ws.Range('E18').Value = '  -5.18%  '
$ws.Range('D19').Value = '39.426.84'
$ws.Range('E19').Value = '  -2.82%  '
$ws.Range('D20').Value = '0.0₃0877'
$ws.Range('E20').Value = '  -3.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.74'
$ws.Range('E21').Value = '  -6.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.27'
$ws.Range('E22').Value = '  -4.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.43'
$ws.Range('E23').Value = '  -2.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '227.49'
$ws.Range('E24').Value = '  -3.14%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  -6.20%  '
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.71'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('E29').Value = '  +0.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.13'
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '149.40'
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.02'
$ws.Range('E32').Value = '  -6.53%  '
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.87'
$ws.Range('E34').Value = '  -6.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0698'
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.110'
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.51'
$ws.Range('E38').Value = '  -2.47%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0970'
$ws.Range('E39').Value = '  -2.64%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.64'
$ws.Range('E40').Value = '  -5.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.65'
$ws.Range('E41').Value = '  -2.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.65'
$ws.Range('E42').Value = '  -4.17%  '
$ws.Range('D43').Value = '1.901.66'
$ws.Range('E43').Value = '  -3.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.06'
$ws.Range('E44').Value = '  -13.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0260'
$ws.Range('E45').Value = '  -2.90%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.24'
$ws.Range('E46').Value = '  -7.61%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.03'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.64'
$ws.Range('E48').Value = '  -2.27%  '
$ws.Range('D49').Value = '2.420.66'
$ws.Range('E49').Value = '  -7.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '70.88'
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '87.47'
$ws.Range('E51').Value = '  -5.91%  '
